# This script updates the division-answer table cells to match the target content.
#
# Quirk of this simulated Word engine: Find.Execute always searches/replaces against
# the WHOLE document content, not just the Range/Cell it was invoked on - and with
# wdReplaceOne (Replace:=1) it always replaces the first (left-most, document-order)
# match, not the match nearest the invoking range. Several of the target table values
# are re-used as the target text of other cells (e.g. '78÷3=26, 0' is both an original
# cell value and the NEW value written into a different cell), so a naive left-to-right
# pass could collide: creating a duplicate value earlier in the document than the real
# occurrence we still need to replace would make wdReplaceOne touch the wrong cell.
#
# To stay correct regardless of collisions we:
#   1. Always replace exactly ONE occurrence at a time (Replace:=1 / wdReplaceOne).
#   2. Order the replacements so a value is never produced as a 'new' value before
#      every pre-existing occurrence of that same text has already been consumed as
#      an 'old' search target (a dependency-safe / topological order).
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("73÷5=14, 3", $true, $false, $false, $false, $false, $true, 1, $false, "63÷6=10, 3", 1)
if (-not $found) { throw "Replace failed for cell R1C1: '73÷5=14, 3' -> '63÷6=10, 3'" }
$found = $d.Content.Find.Execute("30÷5=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷2=41, 0", 1)
if (-not $found) { throw "Replace failed for cell R1C2: '30÷5=6, 0' -> '82÷2=41, 0'" }
$found = $d.Content.Find.Execute("14÷4=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "12÷3=4, 0", 1)
if (-not $found) { throw "Replace failed for cell R1C3: '14÷4=3, 2' -> '12÷3=4, 0'" }
$found = $d.Content.Find.Execute("33÷6=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "61÷6=10, 1", 1)
if (-not $found) { throw "Replace failed for cell R1C4: '33÷6=5, 3' -> '61÷6=10, 1'" }
$found = $d.Content.Find.Execute("59÷2=29, 1", $true, $false, $false, $false, $false, $true, 1, $false, "81÷5=16, 1", 1)
if (-not $found) { throw "Replace failed for cell R1C5: '59÷2=29, 1' -> '81÷5=16, 1'" }
$found = $d.Content.Find.Execute("69÷5=13, 4", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 1)
if (-not $found) { throw "Replace failed for cell R5C1: '69÷5=13, 4' -> '54÷9=6, 0'" }
$found = $d.Content.Find.Execute("22÷3=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=7, 6", 1)
if (-not $found) { throw "Replace failed for cell R5C2: '22÷3=7, 1' -> '69÷9=7, 6'" }
$found = $d.Content.Find.Execute("78÷3=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=15, 2", 1)
if (-not $found) { throw "Replace failed for cell R9C2: '78÷3=26, 0' -> '62÷4=15, 2'" }
$found = $d.Content.Find.Execute("71÷7=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 1)
if (-not $found) { throw "Replace failed for cell R5C3: '71÷7=10, 1' -> '78÷3=26, 0'" }
$found = $d.Content.Find.Execute("64÷9=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "74÷9=8, 2", 1)
if (-not $found) { throw "Replace failed for cell R5C4: '64÷9=7, 1' -> '74÷9=8, 2'" }
$found = $d.Content.Find.Execute("85÷6=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "84÷9=9, 3", 1)
if (-not $found) { throw "Replace failed for cell R5C5: '85÷6=14, 1' -> '84÷9=9, 3'" }
$found = $d.Content.Find.Execute("49÷5=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "99÷9=11, 0", 1)
if (-not $found) { throw "Replace failed for cell R9C1: '49÷5=9, 4' -> '99÷9=11, 0'" }
$found = $d.Content.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "12÷3=4, 0", 1)
if (-not $found) { throw "Replace failed for cell R9C3: '77÷4=19, 1' -> '12÷3=4, 0'" }
$found = $d.Content.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "86÷8=10, 6", 1)
if (-not $found) { throw "Replace failed for cell R9C4: '21÷7=3, 0' -> '86÷8=10, 6'" }
$found = $d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "54÷8=6, 6", 1)
if (-not $found) { throw "Replace failed for cell R9C5: '26÷3=8, 2' -> '54÷8=6, 6'" }
$found = $d.Content.Find.Execute("33÷2=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=1, 4", 1)
if (-not $found) { throw "Replace failed for cell R13C1: '33÷2=16, 1' -> '12÷8=1, 4'" }
$found = $d.Content.Find.Execute("44÷8=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=8, 4", 1)
if (-not $found) { throw "Replace failed for cell R13C2: '44÷8=5, 4' -> '76÷9=8, 4'" }
$found = $d.Content.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "12÷6=2, 0", 1)
if (-not $found) { throw "Replace failed for cell R13C3: '46÷6=7, 4' -> '12÷6=2, 0'" }
$found = $d.Content.Find.Execute("35÷8=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "53÷5=10, 3", 1)
if (-not $found) { throw "Replace failed for cell R13C4: '35÷8=4, 3' -> '53÷5=10, 3'" }
$found = $d.Content.Find.Execute("73÷8=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷8=1, 6", 1)
if (-not $found) { throw "Replace failed for cell R13C5: '73÷8=9, 1' -> '14÷8=1, 6'" }
$found = $d.Content.Find.Execute("31÷9=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "87÷2=43, 1", 1)
if (-not $found) { throw "Replace failed for cell R17C1: '31÷9=3, 4' -> '87÷2=43, 1'" }
$found = $d.Content.Find.Execute("71÷5=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "63÷5=12, 3", 1)
if (-not $found) { throw "Replace failed for cell R17C2: '71÷5=14, 1' -> '63÷5=12, 3'" }
$found = $d.Content.Find.Execute("19÷8=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "98÷4=24, 2", 1)
if (-not $found) { throw "Replace failed for cell R17C3: '19÷8=2, 3' -> '98÷4=24, 2'" }
$found = $d.Content.Find.Execute("66÷8=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "64÷5=12, 4", 1)
if (-not $found) { throw "Replace failed for cell R17C4: '66÷8=8, 2' -> '64÷5=12, 4'" }
$found = $d.Content.Find.Execute("33÷8=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷2=39, 0", 1)
if (-not $found) { throw "Replace failed for cell R17C5: '33÷8=4, 1' -> '78÷2=39, 0'" }
